# Applies the "simple sensor model with observation probabilities" edit:
# updates truth-state observation values in columns B-F of the test
# results sheet, leaving headers, formatting and column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Ship Truth State)
$ws.Range("B6:B22").Value = 2
$ws.Range("B88:B102").Value = 0

# Column C (System 1 Truth State)
$ws.Range("C4:C19").Value = 2

# Column D (System 2 Truth State)
$ws.Range("D6:D22").Value = 2

# Column E (System 3 Truth State)
$ws.Range("E7:E25").Value = 2
$ws.Range("E26:E102").Value = 0

# Column F (System 4 Truth State)
$ws.Range("F6:F34").Value = 2
$ws.Range("F49:F87").Value = 1
